{"js": "// The payment-terms paragraph (\"Sazba za pron\u00e1jem ...\") originally reads:\n//   \" 5 000 K\u010d + DPH, pron\u00e1jem prom\u00edtac\u00edho pl\u00e1tna 1 000 K\u010d + DPH  \"\n// It needs to become:\n//   \" 5 000 K\u010d/den + DPH, pron\u00e1jem prom\u00edtac\u00edho pl\u00e1tna 1 000 K\u010d/den + DPH  \"\n// i.e. insert \"/den\" right after each of the two \"K\u010d\" occurrences that are\n// immediately followed by \" + DPH\".\n\nconst body = context.document.body;\n\n// Find both occurrences of \"K\u010d + DPH\" (there are exactly two, both inside the\n// \"Sazba za pron\u00e1jem...\" paragraph under \u010cl. 3 Platebn\u00ed podm\u00ednky).\nconst matches = body.search(\"K\u010d + DPH\", { matchCase: false, matchWholeWord: false });\nmatches.load(\"items\");\nawait context.sync();\n\n// Walk backwards so inserting text doesn't shift the positions of the other\n// (not-yet-processed) matches.\nfor (let i = matches.items.length - 1; i >= 0; i--) {\n  const match = matches.items[i];\n\n  // Within this match, locate just the \"K\u010d\" part and insert \"/den\" right\n  // after it (formatting is inherited from the surrounding bold run).\n  const kc = match.search(\"K\u010d\", { matchCase: false, matchWholeWord: false });\n  kc.load(\"items\");\n  await context.sync();\n\n  if (kc.items.length > 0) {\n    kc.items[0].insertText(\"/den\", Word.InsertLocation.after);\n  }\n}\n\nawait context.sync();\n", "ps1": "# The payment-terms paragraph (\"Sazba za pron\u00e1jem ...\") originally reads:\n#   \" 5 000 K\u010d + DPH, pron\u00e1jem prom\u00edtac\u00edho pl\u00e1tna 1 000 K\u010d + DPH  \"\n# It needs to become:\n#   \" 5 000 K\u010d/den + DPH, pron\u00e1jem prom\u00edtac\u00edho pl\u00e1tna 1 000 K\u010d/den + DPH  \"\n# i.e. insert \"/den\" right after each of the two \"K\u010d\" occurrences that are\n# immediately followed by \" + DPH\". There are exactly two such occurrences,\n# both inside the \"Sazba za pron\u00e1jem...\" paragraph under \u010cl. 3 Platebn\u00ed podm\u00ednky.\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"K\u010d + DPH\"\n$find.Replacement.Text = \"K\u010d/den + DPH\"\n$find.Forward = $true\n$find.Wrap = 1\n\n# wdReplaceAll = 2 -> replace every occurrence in one pass.\n$find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2)\n"}
